# Insert a new "statut_name" column before the existing NCTId column (column C),
# shifting NCTId..intervention_type one column to the right (C:L -> D:M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()

# Populate the new column's header and data values.
$ws.Range("C1").Value = "statut_name"
$ws.Range("C2:C6").Value = "pas de résultat ni de publication"
